$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.331.03"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.684.94"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.04"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5234"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2697"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06407"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.99"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07498"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "1.712.12"
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.563"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5783"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008472"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.25"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "26.365.38"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.915"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.86"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.51"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.48"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.702"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1231"
$ws.Range("E26").Value = "  +4.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.78"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06625"
$ws.Range("E28").Value = "  +12.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.345"
$ws.Range("E29").Value = "  +6.59%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.569"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.568"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.025"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6211"
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.402"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.365"
$ws.Range("E38").Value = "  +5.70%  "
$ws.Range("D39").Value = "1.108.74"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01614"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8818"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.96"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "1.832.65"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.70"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.173"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4307"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.059"
$ws.Range("E51").Value = "  +3.14%  "
